# really commit from arikan-dell 6/23-2
#
# Applies three logical changes to Sheet1:
#  1. Row 34: the old D34 note ("修改dell") is removed entirely.
#  2. Row 35: a brand-new incident row is filled in (BAIERTE TRADING / Klarna
#     failures) across D:J.
#  3. Rows 46-47: two more brand-new incident rows are filled in
#     (HOPPE TECHNOLOGY CO / HONGXINYI) plus a style-only tweak on C46:C47.
#  4. The sheet view scrolls down and the selection moves to A30:A48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the stray note in D34 (also drops it from shared strings so a
#        later write can reclaim the slot, matching upstream's shared string
#        table). Clear() removes the cell outright (value + style), shrinking
#        row 34's recorded span the same way the source workbook does.
$ws.Range("D34").Clear()

# --- 2. Row 35: new incident - BAIERTE TRADING / Klarna payments all failing.
$ws.Range("D2").Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("D35").Value2 = 45825.4013888889

$ws.Range("E8").Copy()
$ws.Range("E35").PasteSpecial(-4122)
$ws.Range("E35").Value = "BAIERTE TRADING"

$ws.Range("F35").Value = "Klarna付款全部失败"

$ws.Range("G2").Copy()
$ws.Range("G35").PasteSpecial(-4122)
$ws.Range("G35").Value = "用户商户号绑定错误"

$ws.Range("H2").Copy()
$ws.Range("H35").PasteSpecial(-4122)
$ws.Range("H35").Value = "使用最新的商户号和秘钥"

$ws.Range("I3").Copy()
$ws.Range("I35").PasteSpecial(-4122)
$ws.Range("I35").Value2 = 45825.4069444444

$ws.Range("J2").Copy()
$ws.Range("J35").PasteSpecial(-4122)
$ws.Range("J35").Value2 = 45825.4395833333

# --- 3. Row 46: new incident - HOPPE TECHNOLOGY CO checkout issue.
$ws.Range("C3").Copy()
$ws.Range("C46").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("D46").Value2 = 45828.6854166667

$ws.Range("E8").Copy()
$ws.Range("E46").PasteSpecial(-4122)
$ws.Range("E46").Value = "HOPPE TECHNOLOGY CO"

$ws.Range("F46").Value = "客户在网页上直接下单可以正常支付，但是当在后台自己创建的订单时，支付页面没有信息填写选项，无法进入下一步"

$ws.Range("K8").Copy()
$ws.Range("K46").PasteSpecial(-4122)
$ws.Range("K46").Value = "6/20 16:29陈玺接手，6/23 14:28陈玺推送问题到我，6/23 14:34运营刘玲玲通知该问题联系了季超，商户后台给了他测试，于是向季超重新反馈了该问题"

# --- 4. Row 47: new incident - HONGXINYI delayed webhook callback.
$ws.Range("C3").Copy()
$ws.Range("C47").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("D47").Value2 = 45830.6027777778

$ws.Range("E8").Copy()
$ws.Range("E47").PasteSpecial(-4122)
$ws.Range("E47").Value = "HONGXINYI "

$ws.Range("F47").Value = "回调时间太长，隔了几分钟才传过来"

$ws.Range("G2").Copy()
$ws.Range("G47").PasteSpecial(-4122)
$ws.Range("G47").Value = "可能是消息堆积"

$ws.Range("I3").Copy()
$ws.Range("I47").PasteSpecial(-4122)
$ws.Range("I47").Value2 = 45830.6506944444

$ws.Range("K8").Copy()
$ws.Range("K47").PasteSpecial(-4122)
$ws.Range("K47").Value = "东海排查中"

# --- 5. Scroll the sheet view down and select A30:A48, matching the new
#        working area the author was looking at.
$ws.Range("A30:A48").Select()
$excel.ActiveWindow.ScrollRow = 16
